$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 332; existing rows 332:373 shift down to 333:374.
$ws.Rows("332:332").Insert()

# Populate the newly inserted row 332 with the new weekly price record.
$ws.Range("A332").Value = 9
$ws.Range("B332").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C332").Value = "Metropolitana"
$ws.Range("D332").Value = 45142
$ws.Range("E332").Value = 13
$ws.Range("F332").Value = 100112026
$ws.Range("G332").Value = "Haba"
$ws.Range("H332").Value = "Sin especificar"
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 52
$ws.Range("K332").Value = 13000
$ws.Range("L332").Value = 15000
$ws.Range("M332").Value = 14000
$ws.Range("N332").Value = "`$/saco 25 kilos"
$ws.Range("O332").Value = "Provincia de Limarí"
$ws.Range("P332").Value = 560
$ws.Range("Q332").Value = 25
$ws.Range("R332").Value = "Hortaliza"
